# Adds a small two/three-column data table (M11:O18) and a new
# "Soil Moisture Readings (%)" summary table (A28:D31) to the worksheet,
# matching the data added by the commit ("Added a script for running the
# arduino using my laptop").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Small side table in columns M:O (rows 11-18) ---------------------
$ws.Range("M11").Value = "a"
$ws.Range("N11").Value = "b"

$ws.Range("M12").Value = 2
$ws.Range("N12").Value = 123
$ws.Range("O12").Value = 2

$ws.Range("M13").Value = 3
$ws.Range("N13").Value = 23
$ws.Range("O13").Value = 3

$ws.Range("M14").Value = 4
$ws.Range("N14").Value = 12
$ws.Range("O14").Value = 4

$ws.Range("M15").Value = 5
$ws.Range("N15").Value = 43
$ws.Range("O15").Value = 5

$ws.Range("M16").Value = 4
$ws.Range("N16").Value = 5
$ws.Range("O16").Value = 4

$ws.Range("M17").Value = 5
$ws.Range("N17").Value = 6
$ws.Range("O17").Value = 5

$ws.Range("M18").Value = 6
$ws.Range("N18").Value = 7
$ws.Range("O18").Value = 6

# --- New "Soil Moisture Readings (%)" table (A28:D31) ------------------
$ws.Range("A28").Value = "Soil Moisture Readings (%)"
$ws.Range("A28:D28").Merge() | Out-Null

# keep A29 present (blank, default styled) like the target sheet
$ws.Range("A29").HorizontalAlignment = 1

$ws.Range("B29").Value = "Trial 1"
$ws.Range("C29").Value = "Trial 2"
$ws.Range("D29").Value = "Trial 3"

$ws.Range("A30").Value = "Setup A"
$ws.Range("B30").Value = 65

$ws.Range("A31").Value = "Setup B"
$ws.Range("B31").Value = 62

# Bold, centered header styling to match the rest of the workbook's table
# headers (e.g. "Sample Size").
$ws.Range("A28:D28").Font.Bold = $true
$ws.Range("A28:D28").HorizontalAlignment = -4108

$ws.Range("B29:D29").Font.Bold = $true
$ws.Range("B29:D29").HorizontalAlignment = -4108

$ws.Range("A30").Font.Bold = $true
$ws.Range("A30").HorizontalAlignment = -4108

$ws.Range("A31").Font.Bold = $true
$ws.Range("A31").HorizontalAlignment = -4108

# --- Update the active view / selection as in the edited workbook ------
$ws.Range("C30").Select() | Out-Null
